$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the cell content: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the active selection to E8 as shown in the updated sheetView
$ws.Activate()
$ws.Range("E8").Select()
